# "Generate Report for Handback"
#
# For each locale sheet (zh-cn, de-de) the handback pipeline has now run:
#   - Status (column C) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every data row.
#   - Latest Target File (I) / Latest Handback File (J) / Latest Handback
#     DateTime (K) get populated for every data row.
#   - The Latest Target File cell becomes a hyperlink to the source file
#     (mirroring the existing Source File Name hyperlink), with display
#     text "a.md".
#   - A couple of report columns get widened so the new long file names /
#     status text fit.

function Update-LocaleSheet($wb, $repoBlobBase, $SheetName, $HandoffXlf, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the Status column (C) and the Latest Handback File column (J)
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Report rows 2 and 3 (a.md / b.md) are now handed back and in sync.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Handback File (J) and Latest Handback DateTime (K).
    $ws.Range("J2").Value = $HandoffXlf
    $ws.Range("J3").Value = $HandoffXlf
    $ws.Range("K2").Value = $HandbackDateTime
    $ws.Range("K3").Value = $HandbackDateTime

    # Latest Target File (I) becomes a hyperlink back to the source file,
    # same target/display as the Source File Name hyperlink in column A.
    # Rebuild the hyperlink collection so the links land in document order
    # (A2, I2, A3, I3).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), ($repoBlobBase + "a.md"), "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), ($repoBlobBase + "a.md"), "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), ($repoBlobBase + "b.md"), "", "", "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), ($repoBlobBase + "a.md"), "", "", "a.md")
}

$wb = $excel.ActiveWorkbook

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d0fd6d22951f2572bd4cfcf520aec50f607d8d1/e2e/"

# --- Overview sheet: widen the zh-cn / de-de status columns (E, F) -------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn handed back first.
Update-LocaleSheet $wb $repoBlobBase "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-09-04 20:43:19"

# de-de handed back a few seconds later.
Update-LocaleSheet $wb $repoBlobBase "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-09-04 20:43:26"
